$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A37").Value = 45883
Write-Host "done"
